$d = $word.ActiveDocument

# Locate the paragraph that ends with "...(Requisito fraco)" — the block to
# remove (a blank paragraph, the "Ver no Jupiter..." paragraph, and the
# "© 2020 ..." footer paragraph) starts right after it.
$startRange = $d.Content
$startRange.Find.Execute("LOQ4240: Administração e Organização II (Requisito fraco)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startRange.MoveEnd(1, 1) | Out-Null   # include the paragraph mark
$startDelete = $startRange.End

# Locate the end of the footer paragraph ("© 2020 ... Creative Commons
# Attribution"), which is the last paragraph to be removed.
$endRange = $d.Content
$endRange.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endRange.MoveEnd(1, 1) | Out-Null     # include the paragraph mark
$endDelete = $endRange.End

# Remove the blank paragraph + the two trailing paragraphs (the leftover
# empty paragraph that used to follow them is left in place).
$d.Range($startDelete, $endDelete).Delete()
